$wb = $excel.ActiveWorkbook

# --- 1. Insert a new "Player Info" sheet before the current first sheet ---
$firstSheet = $wb.Worksheets.Item(1)
$playerInfo = $wb.Worksheets.Add($firstSheet)
$playerInfo.Name = "Player Info"

# Re-fetch the existing sheets by name now that the sheet collection has
# shifted (sheet references captured before the insert point at a stale
# position rather than following the renamed/moved sheet).
$battingSheet = $wb.Worksheets.Item("ODI Batting")
$bowlingSheet = $wb.Worksheets.Item("ODI Bowling")

# --- 2. Populate "Player Info" ---
$playerInfo.Range("A1").Value = "ID"
$playerInfo.Range("B1").Value = "NAME"
$playerInfo.Range("C1").Value = "BATTING_HAND"
$playerInfo.Range("D1").Value = "BOWL_STYLE"
# Reuse the same header style as the other sheets (bold, centered, bordered)
# by copying the format from an existing header cell instead of re-deriving
# the same look via discrete property sets (which would mint a new,
# merely-equivalent style entry).
$battingSheet.Range("A1").Copy()
$playerInfo.Range("A1:D1").PasteSpecial(-4122)

$playerInfo.Range("A2").Value = "'5839"
$playerInfo.Range("A2").Style = "Normal"
$playerInfo.Range("B2").Value = "Daryn Miles Dupavillon"
$playerInfo.Range("C2").Value = "Right Handed"
$playerInfo.Range("D2").Value = "Right Arm Fast"

# --- 3. ODI Batting: MATCH_CARD_LINK -> MATCH_CODE, url -> bare match code ---
$battingSheet.Range("D1").Value = "MATCH_CODE"
$battingSheet.Range("D2").Value = "'4421"
$battingSheet.Range("D2").Style = "Normal"
$battingSheet.Range("D3").Value = "'4460"
$battingSheet.Range("D3").Style = "Normal"

# --- 4. ODI Bowling: MATCH_CARD_LINK -> MATCH_CODE, url -> bare match code ---
$bowlingSheet.Range("B1").Value = "MATCH_CODE"
$bowlingSheet.Range("B2").Value = "'4421"
$bowlingSheet.Range("B2").Style = "Normal"
$bowlingSheet.Range("B3").Value = "'4460"
$bowlingSheet.Range("B3").Style = "Normal"
